# Auto-generated edit script: updates Leve profit-calculation values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (per the scheduled-runner data refresh).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1016.6923
$ws.Range("I41").Value = 1264.2858
$ws.Range("J41").Value = 727.8333
$ws.Range("K41").Value = 1264.2858
$ws.Range("L41").Value = 727.8333
$ws.Range("M41").Value = -824.2858000000001
$ws.Range("N41").Value = -1607.8333
$ws.Range("H138").Value = 2881.1223
$ws.Range("I138").Value = 3507.1
$ws.Range("J138").Value = 2802.875
$ws.Range("K138").Value = 10521.3
$ws.Range("L138").Value = 8408.625
$ws.Range("M138").Value = -5381.299999999999
$ws.Range("N138").Value = -18688.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24635.807
$ws.Range("I32").Value = 7432.976
$ws.Range("K32").Value = 7432.976
$ws.Range("M32").Value = -7145.976
$ws.Range("H74").Value = 826.41174
$ws.Range("I74").Value = 621.4545000000001
$ws.Range("K74").Value = 621.4545000000001
$ws.Range("M74").Value = 252.5454999999999
$ws.Range("H77").Value = 826.41174
$ws.Range("I77").Value = 621.4545000000001
$ws.Range("K77").Value = 3107.2725
$ws.Range("M77").Value = 1260.7275
$ws.Range("H122").Value = 1960.4572
$ws.Range("I122").Value = 1648.5416
$ws.Range("K122").Value = 4945.6248
$ws.Range("M122").Value = -2495.6248
$ws.Range("H132").Value = 12765.755
$ws.Range("I132").Value = 16440.975
$ws.Range("J132").Value = 2527.6428
$ws.Range("K132").Value = 49322.925
$ws.Range("L132").Value = 7582.928400000001
$ws.Range("M132").Value = -46792.925
$ws.Range("N132").Value = -12642.9284

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2035.4286
$ws.Range("I99").Value = 1187.7778
$ws.Range("J99").Value = 2328.8462
$ws.Range("K99").Value = 1187.7778
$ws.Range("L99").Value = 2328.8462
$ws.Range("M99").Value = 310.2221999999999
$ws.Range("N99").Value = -5324.8462
$ws.Range("H134").Value = 2692.5918
$ws.Range("I134").Value = 2622.75
$ws.Range("K134").Value = 7868.25
$ws.Range("M134").Value = -5333.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32552.064
$ws.Range("J31").Value = 58063.152
$ws.Range("L31").Value = 58063.152
$ws.Range("N31").Value = -58653.152
$ws.Range("H34").Value = 32552.064
$ws.Range("J34").Value = 58063.152
$ws.Range("L34").Value = 58063.152
$ws.Range("N34").Value = -58467.152
$ws.Range("H107").Value = 3846.6875
$ws.Range("I107").Value = 5244.591
$ws.Range("J107").Value = 771.3
$ws.Range("K107").Value = 5244.591
$ws.Range("L107").Value = 771.3
$ws.Range("M107").Value = -3324.591
$ws.Range("N107").Value = -4611.3
$ws.Range("H132").Value = 3699.2778
$ws.Range("I132").Value = 3536.8125
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 10610.4375
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -8080.4375
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 1742.8572
$ws.Range("I134").Value = 1325
$ws.Range("J134").Value = 2300
$ws.Range("K134").Value = 3975
$ws.Range("L134").Value = 6900
$ws.Range("M134").Value = -1440
$ws.Range("N134").Value = -11970

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 600
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 1800
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -2146
$ws.Range("H34").Value = 1064.8
$ws.Range("J34").Value = 1324.75
$ws.Range("L34").Value = 3974.25
$ws.Range("N34").Value = -4142.25
$ws.Range("H39").Value = 2674.6667
$ws.Range("J39").Value = 2946.5
$ws.Range("L39").Value = 8839.5
$ws.Range("N39").Value = -9427.5
$ws.Range("H68").Value = 1072.1666
$ws.Range("I68").Value = 607.5
$ws.Range("J68").Value = 2001.5
$ws.Range("K68").Value = 1822.5
$ws.Range("L68").Value = 6004.5
$ws.Range("M68").Value = -1011.5
$ws.Range("N68").Value = -7626.5
$ws.Range("H71").Value = 1072.1666
$ws.Range("I71").Value = 607.5
$ws.Range("J71").Value = 2001.5
$ws.Range("K71").Value = 5467.5
$ws.Range("L71").Value = 18013.5
$ws.Range("M71").Value = -1411.5
$ws.Range("N71").Value = -26125.5
$ws.Range("H106").Value = 2401.3333
$ws.Range("J106").Value = 2401.3333
$ws.Range("L106").Value = 7203.999899999999
$ws.Range("N106").Value = -9095.999899999999
$ws.Range("H108").Value = 1819.7
$ws.Range("I108").Value = 1371.875
$ws.Range("J108").Value = 3611
$ws.Range("K108").Value = 4115.625
$ws.Range("L108").Value = 10833
$ws.Range("M108").Value = -1235.625
$ws.Range("N108").Value = -16593
$ws.Range("H109").Value = 3825.8
$ws.Range("I109").Value = 10000
$ws.Range("J109").Value = 3568.5417
$ws.Range("K109").Value = 30000
$ws.Range("L109").Value = 10705.6251
$ws.Range("M109").Value = -28960
$ws.Range("N109").Value = -12785.6251
$ws.Range("H130").Value = 800
$ws.Range("I130").Value = 600
$ws.Range("J130").Value = 1000
$ws.Range("K130").Value = 1800
$ws.Range("L130").Value = 3000
$ws.Range("M130").Value = 3220
$ws.Range("N130").Value = -13040
$ws.Range("H131").Value = 649772.9399999999
$ws.Range("I131").Value = 648
$ws.Range("J131").Value = 745232.5
$ws.Range("K131").Value = 1944
$ws.Range("L131").Value = 2235697.5
$ws.Range("M131").Value = 3096
$ws.Range("N131").Value = -2245777.5
$ws.Range("H137").Value = 3613231.5
$ws.Range("J137").Value = 5890939.5
$ws.Range("L137").Value = 17672818.5
$ws.Range("N137").Value = -17683018.5
$ws.Range("H139").Value = 2118.125
$ws.Range("I139").Value = 1130
$ws.Range("J139").Value = 2886.6667
$ws.Range("K139").Value = 3390
$ws.Range("L139").Value = 8660.000100000001
$ws.Range("M139").Value = 1750
$ws.Range("N139").Value = -18940.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 121379.86
$ws.Range("I102").Value = 1136.3784
$ws.Range("J102").Value = 463611.3
$ws.Range("K102").Value = 1136.3784
$ws.Range("L102").Value = 463611.3
$ws.Range("M102").Value = 485.6215999999999
$ws.Range("N102").Value = -466855.3

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 100007
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 100007
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100007
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -100287
$ws.Range("H136").Value = 1885.85
$ws.Range("I136").Value = 1289.5
$ws.Range("K136").Value = 3868.5
$ws.Range("M136").Value = -1318.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 166680340
$ws.Range("I39").Value = 1000000000
$ws.Range("J39").Value = 16400
$ws.Range("K39").Value = 1000000000
$ws.Range("L39").Value = 16400
$ws.Range("M39").Value = -999999587
$ws.Range("N39").Value = -17226
$ws.Range("H56").Value = 25053.334
$ws.Range("I56").Value = 7463.3335
$ws.Range("K56").Value = 7463.3335
$ws.Range("M56").Value = -6749.3335
$ws.Range("H136").Value = 14799.915
$ws.Range("I136").Value = 24410.785
$ws.Range("J136").Value = 4708.5
$ws.Range("K136").Value = 73232.355
$ws.Range("L136").Value = 14125.5
$ws.Range("M136").Value = -70682.355
$ws.Range("N136").Value = -19225.5

Write-Host "Applied 183 value updates and 2 clears."
